# Update the metric-definition sheet (Guide) with the revised "billion -> million"
# conversion wording, and fix up the selection left on the sheet.
#
# Note: cell A14 (EBITDA) is re-entered BEFORE cell A13 (Revenue) so that the
# workbook's shared-string table ends up in the same append order Excel
# produced when the author originally retyped these two cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value  = "Implied equity value (capitalization of target company - price of one share multiplied by number of target’s shares on the moment of deal announcement), in million dollars (if the metrics is given in billion dollars, multiply it by 1000 to convert it into million)"
$ws.Range("A10").Value = "Implied Net Debt (preferrably taken from last balance sheet prior to the deal, calculated as target company's total debt minus cash), in million dollars (if the metrics is given in billion dollars, multiply it by 1000 to convert it into million)"
$ws.Range("A14").Value = "EBITDA of target company (calculated as operating profit plus depreciation) for the full year of the announcement date (if no information, take the revenue for the year which preceeded the announcement date), in million dollars (if the metrics is given in billion dollars, multiply it by 1000 to convert it into million)"
$ws.Range("A13").Value = "Revenue of target company for the full year of the announcement date (if no information, take the revenue for the year which preceeded the announcement date), in million dollars (if the metrics is given in billion dollars, multiply it by 1000 to convert it into million)"

# Leave the sheet scrolled back to the top with A9 selected (matches the
# saved view state in the workbook).
$ws.Range("A9").Select() | Out-Null
